# Add "NA" values under the duplicate_image_filename column (column E)
# for the practice/trial rows 2-21 of the stimuli table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

foreach ($r in 2..21) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
